# Remove the block of trailing empty paragraphs that follows the
# "...tutorial!" paragraph, keeping only the very last (already empty)
# paragraph that precedes the section properties.
#
# Before:
#   ...tutorial!
#   <empty ListParagraph, b/bCs sz=36>
#   <empty, sz=24>
#   <empty, b/bCs sz=40>
#   <empty>
#   <empty>
#   <empty>
#   <empty>                <- kept
#   sectPr
#
# After:
#   ...tutorial!
#   <empty>                <- kept
#   sectPr

$d = $word.ActiveDocument

# Locate the paragraph that ends the "...tutorial!" sentence; it is the
# anchor right before the run of empty paragraphs that must be collapsed.
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*tutorial!*") {
        $anchorIndex = $i
    }
}

if ($anchorIndex -gt 0) {
    $total = $d.Paragraphs.Count

    # Keep the final paragraph in the document (it precedes sectPr); delete
    # everything between the anchor paragraph and that last paragraph.
    if ($total -gt ($anchorIndex + 1)) {
        $startPara = $d.Paragraphs.Item($anchorIndex + 1)
        $endPara = $d.Paragraphs.Item($total - 1)

        $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
        $r.Delete()
    }
}
